# Predicted_Ekstraklasa2025_26_table_matchday_5.xlsx
#
# Insert four new columns (C:F) before the existing "ExpPoints" column
# (which shifts from C to G) and label them WIN / TOP2 / TOP4 / RELEGATION.
# These are placeholders for a future Monte Carlo simulation (per the
# commit message) - for now they are left blank for every team row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank columns at C:F - this shifts the old column C (ExpPoints)
# to column G and carries the header's style (bold/bordered/centered) along
# into the newly inserted header cells.
$ws.Columns("C:F").Insert()

# New header labels for the inserted columns.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"

# Leave the data rows for the new columns blank (no values yet - these will
# be filled in later by the Monte Carlo simulation).
$ws.Range("C2:F19").Value = ""
